# Edit: for each year block of 4 rows (A/B/C/D quarter rows), swap the
# "B" quarter row and "C" quarter row (columns A:E), then drop columns F:G
# (the "产销率"/"销售量" duplicate columns) entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 81
$rowsPerBlock = 4
$lastCol = 5   # columns A..E

for ($baseRow = $firstDataRow; $baseRow -le $lastDataRow; $baseRow += $rowsPerBlock) {
    $rowB = $baseRow + 1
    $rowC = $baseRow + 2

    $valsB = @()
    $valsC = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $valsB += , $ws.Cells.Item($rowB, $c).Value2
        $valsC += , $ws.Cells.Item($rowC, $c).Value2
    }

    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($rowB, $c).Value = $valsC[$c - 1]
        $ws.Cells.Item($rowC, $c).Value = $valsB[$c - 1]
    }
}

# Drop the redundant F (产销率) and G (销售量) columns entirely.
$ws.Range("F1:G1").EntireColumn.Delete()
